$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B43").Value = 664
$ws.Range("B44").Value = 657
$ws.Range("B45").Value = 691
$ws.Range("B46").Value = 769
$ws.Range("B47").Value = 829
$ws.Range("B48").Value = 845
$ws.Range("B49").Value = 819
$ws.Range("B50").Value = 797
